$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$row = 56
$logs.Cells.Item($row, 1).Value = "Inlogproblemen"
$logs.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($row, 3).Value = "Ik kan niet meer inloggen op mijn account. Kunnen jullie helpen?"
$logs.Cells.Item($row, 4).Value = "IT / Technisch probleem"
$logs.Cells.Item($row, 5).Value = "Beste klant,`nBedankt voor je bericht. Om je verder te kunnen helpen met het inlogprobleem, hebben we wat meer informatie nodig. Zou je alsjeblieft je gebruikersnaam en eventuele foutmeldingen die je hebt ontvangen kunnen delen? Op die manier kunnen we het probleem zo snel mogelijk voor je oplossen.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Cells.Item($row, 6).Value = "2025-06-22 22:09:26"
$logs.Cells.Item($row, 7).Value = "Ja"

# The multi-line answer in column E otherwise triggers an auto row-height
# bump on the newly inserted row; AutoFit drops the resulting customHeight
# override so the row height stays the sheet default, matching the source row.
$logs.Rows.Item($row).AutoFit()

$dashboard.Cells.Item(2, 2).Value = 10

# Extend the conditional-formatting ranges to cover the new row.
$dFc = $logs.Range("D2:D55").FormatConditions
$dFc.Item(1).ModifyAppliesToRange($logs.Range("D2:D56"))

$gFc = $logs.Range("G2:G55").FormatConditions
$gFc.Item(1).ModifyAppliesToRange($logs.Range("G2:G56"))

